$d = $word.ActiveDocument

# --- Step 1 ---------------------------------------------------------------
# Split the final paragraph right after "...hjälpa William." so the old
# paragraph mark (which carries the _GoBack bookmark) becomes its own,
# now-empty, trailing paragraph.
$r1 = $d.Content
$null = $r1.Find.Execute("hjälpa William.", $false, $false, $false, $false, `
    $false, $true, 1, $false, "hjälpa William.^p", 2)

# --- Step 2 ---------------------------------------------------------------
# Append a new paragraph at the very end of the document and give it the
# "24-okt" text.
$r2 = $d.Content
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$rr2 = $p2.Range
$null = $rr2.MoveEnd(1, -1)
$rr2.Text = "24-okt"

# --- Step 3 ---------------------------------------------------------------
# Append the "I fredags: ..." paragraph.
$r3 = $d.Content
$r3.Collapse(0)
$r3.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$rr3 = $p3.Range
$null = $rr3.MoveEnd(1, -1)
$rr3.Text = "I fredags:  Emma  och Dennis fixade ett förmulär för att ge " + `
    "rättigheter samt snyggade till sidan. Emil och Pontus gjorde klart " + `
    "redigering av poster. William blev klar med sessions för " + `
    "kommentarerna.  Hannes började med felsökning. "

# --- Step 4 ---------------------------------------------------------------
# Append the final "Idag: Dennis ..." paragraph. A placeholder character
# ("X") is appended at the very end so the _GoBack bookmark can be
# (re)created just before it without hitting the edge case where adding a
# zero-length bookmark exactly at the document's end resets it; the
# placeholder is then deleted, leaving the bookmark collapsed in the
# correct spot.
$r4 = $d.Content
$r4.Collapse(0)
$r4.InsertParagraphAfter()
$p4 = $d.Paragraphs.Item($d.Paragraphs.Count)
$rr4 = $p4.Range
$null = $rr4.MoveEnd(1, -1)
$rr4.Text = "Idag: Dennis ska fixa registrerings formuläret igen för att " + `
    "stödja ÅÄÖ, William ska bli klar med kommentarerna. William ska " + `
    "fortsätta med felsökning. Emma blir ledare tills på onsdag kl 12 " + `
    "och fixa så att man ser bloggar även om man inte är inloggad samt " + `
    "ta bort rättigheter. Emil ska fixa en follower counter.  Pontus " + `
    "ska fixa edit sakerna.X"

# --- Step 5 ---------------------------------------------------------------
# Move the _GoBack bookmark from the old (now mid-document) paragraph to
# the end of the new last paragraph, right before its paragraph mark.
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$rrLast = $pLast.Range
$null = $rrLast.MoveEnd(1, -1)
$lastCharStart = $rrLast.End - 1
$bmRange = $d.Range($lastCharStart, $rrLast.End)

$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the placeholder character now that the bookmark anchors on it.
$delRange = $d.Range($lastCharStart, $lastCharStart + 1)
$delRange.Delete()
